# Update trading signals with new market data and signal statuses

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Active Signals
# ---------------------------------------------------------------------------
$active = $wb.Worksheets.Item("Active Signals")

# Row 2: XAUAUD BUY - refreshed market data
$active.Range("A2").Value = "2025-07-28 21:28"
$active.Range("B2").Value = "XAUAUD"
$active.Range("C2").Value = "BUY"
$active.Range("D2").Value = 4064.91481
$active.Range("E2").Value = 4064.91121
$active.Range("F2").Value = 4064.91896
$active.Range("G2").Value = 0.01
$active.Range("H2").Value = "84.0%"
$active.Range("I2").Value = 1.15
$active.Range("J2").Value = "Active"

# Row 3: now XAUGBP SELL (was XAUCAD BUY) - recolor Signal cell to the
# "SELL" red/pink fill used elsewhere in the workbook
$active.Range("A3").Value = "2025-07-28 20:55"
$active.Range("B3").Value = "XAUGBP"
$active.Range("C3").Value = "SELL"
$active.Range("C3").Interior.Color = 13551615
$active.Range("D3").Value = 2109.70362
$active.Range("E3").Value = 2109.7061
$active.Range("F3").Value = 2109.69605
$active.Range("G3").Value = 0.09
$active.Range("H3").Value = "84.0%"
$active.Range("I3").Value = 3.04
$active.Range("J3").Value = "Active"

# Rows 4-7 no longer exist in the refreshed signal set
$active.Rows.Item(4).Resize(4).Delete()

# ---------------------------------------------------------------------------
# Sheet: Summary Dashboard
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary Dashboard")

$summary.Range("B4").Value = 2
$summary.Range("B5").Value = 7
$summary.Range("B6").Value = 8
$summary.Range("B7").Value = "83.5%"
$summary.Range("B8").Value = "1.93"
$summary.Range("B9").Value = "2025-07-28 21:07:35"

# ---------------------------------------------------------------------------
# Sheet: Signal History
# ---------------------------------------------------------------------------
$history = $wb.Worksheets.Item("Signal History")

$historyRows = @(
    @{ Row=2;  A="2025-07-28 21:19"; B="XAUUSD"; C="BUY";  D=2638.81797; E=2638.81536; F=2638.82472; G=0.03; H=0.91; I=2.59; J="Filled" },
    @{ Row=3;  A="2025-07-28 20:56"; B="NZDUSD"; C="SELL"; D=0.5863;     E=0.58862;    F=0.58136;    G=0.1;  H=0.85; I=2.12; J="Pending" },
    @{ Row=4;  A="2025-07-28 20:43"; B="EURUSD"; C="SELL"; D=1.10395;    E=1.10659;    F=1.09987;    G=0.02; H=0.78; I=1.54; J="Filled" },
    @{ Row=5;  A="2025-07-28 20:54"; B="XAUCAD"; C="SELL"; D=3602.6381;  E=3602.64162; F=3602.63223; G=0.1;  H=0.85; I=1.67; J="Pending" },
    @{ Row=6;  A="2025-07-28 21:28"; B="XAUAUD"; C="BUY";  D=4064.91481; E=4064.91121; F=4064.91896; G=0.01; H=0.84; I=1.15; J="Active" },
    @{ Row=7;  A="2025-07-28 20:55"; B="XAUGBP"; C="SELL"; D=2109.70362; E=2109.7061;  F=2109.69605; G=0.09; H=0.84; I=3.04; J="Active" },
    @{ Row=8;  A="2025-07-28 20:40"; B="XAUCAD"; C="BUY";  D=3637.04486; E=3637.0413;  F=3637.05461; G=0.09; H=0.77; I=2.75; J="Pending" },
    @{ Row=9;  A="2025-07-28 21:33"; B="XAUCHF"; C="BUY";  D=2330.19431; E=2330.19843; F=2330.18961; G=0.09; H=0.79; I=1.14; J="Filled" },
    @{ Row=10; A="2025-07-28 21:19"; B="XAUUSD"; C="BUY";  D=2654.13881; E=2654.13442; F=2654.14534; G=0.06; H=0.75; I=1.48; J="Pending" },
    @{ Row=11; A="2025-07-28 20:52"; B="XAUEUR"; C="BUY";  D=2414.83832; E=2414.84059; F=2414.82938; G=0.03; H=0.77; I=3.94; J="Filled" },
    @{ Row=12; A="2025-07-28 20:59"; B="NZDUSD"; C="BUY";  D=0.58938;    E=0.59428;    F=0.58413;    G=0.02; H=0.89; I=1.07; J="Pending" },
    @{ Row=13; A="2025-07-28 21:25"; B="EURUSD"; C="BUY";  D=1.10743;    E=1.10362;    F=1.1122;     G=0.04; H=0.77; I=1.25; J="Pending" },
    @{ Row=14; A="2025-07-28 21:30"; B="XAUEUR"; C="BUY";  D=2412.942;   E=2412.93763; F=2412.94933; G=0.07000000000000001; H=0.9;  I=1.68; J="Filled" },
    @{ Row=15; A="2025-07-28 21:35"; B="USDCAD"; C="BUY";  D=1.3615;     E=1.36633;    F=1.35425;    G=0.02; H=0.87; I=1.5;  J="Pending" },
    @{ Row=16; A="2025-07-28 21:07"; B="USDCHF"; C="BUY";  D=0.88156;    E=0.87934;    F=0.88597;    G=0.03; H=0.9399999999999999; I=1.99; J="Filled" }
)

foreach ($r in $historyRows) {
    $row = $r.Row
    $history.Range("A$row").Value = $r.A
    $history.Range("B$row").Value = $r.B
    $history.Range("C$row").Value = $r.C
    $history.Range("D$row").Value = $r.D
    $history.Range("E$row").Value = $r.E
    $history.Range("F$row").Value = $r.F
    $history.Range("G$row").Value = $r.G
    $history.Range("H$row").Value = $r.H
    $history.Range("I$row").Value = $r.I
    $history.Range("J$row").Value = $r.J
}
